# edit.ps1
# Applies feature-selection table updates to both worksheets
# ("final_fail" and "final_gifted") per the commit diff: re-derived
# boolean selection flags, updated Total counts, and corrected/
# reordered Feature names.

$wb = $excel.ActiveWorkbook

$wsFail = $wb.Worksheets.Item("final_fail")
$wsGifted = $wb.Worksheets.Item("final_gifted")

$dataFail = @(
    ,('Average grade of assignments', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Largest period of inactivity (h)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Submissions (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Clicks (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('On/off campus click ratio', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Days with no interaction', $true, $true, $false, $true, $true, $false, $true, $false, 5)
    ,('Start of Session 1 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5)
    ,('Clicks per session', $true, $true, $false, $true, $true, $false, $true, $false, 5)
    ,('Number of days', $true, $true, $false, $true, $true, $false, $true, $false, 5)
    ,('Resources viewed', $true, $true, $true, $false, $true, $false, $true, $false, 5)
    ,('Number of clicks', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Clicks on folder', $true, $true, $true, $false, $false, $false, $true, $false, 4)
    ,('Start of Session 3 (%)', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Clicks on campus', $false, $false, $true, $true, $true, $false, $true, $false, 4)
    ,('Total time online (min)', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Average session duration (min)', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Start of Session 2 (%)', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Start of Session 7 (%)', $true, $false, $false, $false, $true, $false, $true, $false, 3)
    ,('Clicks per day', $false, $false, $true, $false, $true, $false, $true, $false, 3)
    ,('Links viewed', $true, $true, $false, $false, $false, $false, $true, $false, 3)
    ,('Assignments submitted', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Assignments viewed', $false, $true, $false, $false, $false, $false, $true, $false, 2)
    ,('Discussions viewed', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Forum posts', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Files downloaded', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Start of Session 4 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Number of sessions', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Days with no interaction (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Quizzes started', $false, $false, $false, $true, $false, $false, $true, $false, 2)
    ,('Start of Session 5 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Start of Session 8 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Clicks on course', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Clicks on forum', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Start of Session 6 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Start of Session 10 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Start of Session 9 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
)

$dataGifted = @(
    ,('Average grade of assignments', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Largest period of inactivity (h)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Total time online (min)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Average session duration (min)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Start of Session 1 (%)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Clicks (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6)
    ,('Resources viewed', $true, $true, $false, $true, $true, $false, $true, $false, 5)
    ,('Clicks per session', $true, $true, $true, $false, $true, $false, $true, $false, 5)
    ,('On/off campus click ratio', $true, $true, $false, $true, $true, $false, $true, $false, 5)
    ,('Days with no interaction', $true, $true, $false, $true, $true, $false, $true, $false, 5)
    ,('Assignments submitted', $true, $true, $true, $false, $false, $false, $true, $false, 4)
    ,('Number of days', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Clicks per day', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Start of Session 6 (%)', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Assignments viewed', $true, $true, $false, $false, $true, $false, $true, $false, 4)
    ,('Start of Session 4 (%)', $true, $false, $false, $false, $true, $false, $true, $false, 3)
    ,('Days with no interaction (%)', $true, $false, $false, $false, $true, $false, $true, $false, 3)
    ,('Start of Session 2 (%)', $false, $true, $false, $false, $true, $false, $true, $false, 3)
    ,('Submissions (% of course total)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Files downloaded', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Quizzes started', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Number of clicks', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Clicks on folder', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Clicks on course', $true, $false, $false, $false, $false, $false, $true, $false, 2)
    ,('Start of Session 7 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Start of Session 5 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Start of Session 3 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Clicks on campus', $false, $false, $false, $false, $true, $false, $true, $false, 2)
    ,('Links viewed', $false, $false, $true, $false, $false, $false, $true, $false, 2)
    ,('Discussions viewed', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Forum posts', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Number of sessions', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Clicks on forum', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Start of Session 10 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Start of Session 9 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
    ,('Start of Session 8 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1)
)


$r = 2
foreach ($row in $dataFail) {
    $wsFail.Cells.Item($r, 1).Value = $row[0]
    $wsFail.Cells.Item($r, 2).Value = $row[1]
    $wsFail.Cells.Item($r, 3).Value = $row[2]
    $wsFail.Cells.Item($r, 4).Value = $row[3]
    $wsFail.Cells.Item($r, 5).Value = $row[4]
    $wsFail.Cells.Item($r, 6).Value = $row[5]
    $wsFail.Cells.Item($r, 7).Value = $row[6]
    $wsFail.Cells.Item($r, 8).Value = $row[7]
    $wsFail.Cells.Item($r, 9).Value = $row[8]
    $wsFail.Cells.Item($r, 10).Value = $row[9]
    $r = $r + 1
}

$r = 2
foreach ($row in $dataGifted) {
    $wsGifted.Cells.Item($r, 1).Value = $row[0]
    $wsGifted.Cells.Item($r, 2).Value = $row[1]
    $wsGifted.Cells.Item($r, 3).Value = $row[2]
    $wsGifted.Cells.Item($r, 4).Value = $row[3]
    $wsGifted.Cells.Item($r, 5).Value = $row[4]
    $wsGifted.Cells.Item($r, 6).Value = $row[5]
    $wsGifted.Cells.Item($r, 7).Value = $row[6]
    $wsGifted.Cells.Item($r, 8).Value = $row[7]
    $wsGifted.Cells.Item($r, 9).Value = $row[8]
    $wsGifted.Cells.Item($r, 10).Value = $row[9]
    $r = $r + 1
}
